$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'65.922.98"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.75%  '
$ws.Range('D3').Value = "'3.479.73"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'582.86"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').Value = "'174.00"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = "'0.596"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.85%  '
$ws.Range('D9').Value = "'3.479.94"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('E10').Value = '  -5.42%  '
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('E12').Value = '  -3.72%  '
$ws.Range('D13').Value = "'4.083.76"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').Value = "'30.14"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.70%  '
$ws.Range('D16').Value = "'66.026.64"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = "'3.481.54"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E19').Value = '  -3.62%  '
$ws.Range('D20').Value = "'13.99"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = "'366.57"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.25%  '
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = "'72.54"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('D25').Value = "'0.537"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  +5.11%  '
$ws.Range('D27').Value = "'9.64"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.92%  '
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').Value = "'24.09"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('D31').Value = "'5.78"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = "'7.15"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('E35').Value = '  -7.08%  '
$ws.Range('D36').Value = "'1.55"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').Value = "'160.05"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('D38').Value = "'29.23"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.74%  '
$ws.Range('D39').Value = "'0.889"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').Value = "'2.831.25"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.94%  '
$ws.Range('E41').Value = '  -4.97%  '
$ws.Range('D42').Value = "'2.61"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.42%  '
$ws.Range('D43').Value = "'4.46"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.61%  '
$ws.Range('D44').Value = "'6.44"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('D46').Value = "'39.92"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.27%  '
$ws.Range('D47').Value = "'24.19"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.29%  '
$ws.Range('D48').Value = "'0.0289"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').Value = "'310.64"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.58%  '
$ws.Range('D50').Value = "'0.823"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('E51').Value = '  -1.80%  '
